$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for columns B:G, rows 2:11
$oldValues = @{}
for ($r = 2; $r -le 11; $r++) {
    $rowVals = @()
    for ($c = 2; $c -le 7; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $oldValues[$r] = $rowVals
}

# Shift rows down: new row r (3..11) gets old row (r-1) values
for ($r = 11; $r -ge 3; $r--) {
    $src = $oldValues[$r - 1]
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $src[$c - 2]
    }
}

# Set brand new values for row 2
$ws.Cells.Item(2, 2).Value = -0.09914768942982544
$ws.Cells.Item(2, 3).Value = 0.6760636891099804
$ws.Cells.Item(2, 4).Value = 1.058759192261776
$ws.Cells.Item(2, 5).Value = 1.02896024814459
$ws.Cells.Item(2, 6).Value = 1.052237025233689
$ws.Cells.Item(2, 7).Value = 19
